$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5801
$ws.Range("J3").Value = 6190
$ws.Range("J4").Value = 1343
$ws.Range("J5").Value = 473
$ws.Range("J6").Value = 7904
$ws.Range("J7").Value = 21711

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 417
$ws.Range("J6").Value = 461
$ws.Range("J7").Value = 1365

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 133
$ws.Range("J7").Value = 440

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 242
$ws.Range("J3").Value = 336
$ws.Range("J6").Value = 345
$ws.Range("J7").Value = 1007

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 110
$ws.Range("J7").Value = 320

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 224
$ws.Range("J7").Value = 667

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 161
$ws.Range("J7").Value = 550

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J7").Value = 639
$ws.Range("J8").Value = 1365
$ws.Range("J9").Value = 106
$ws.Range("J11").Value = 340
$ws.Range("J15").Value = 240
$ws.Range("J19").Value = 639
$ws.Range("J20").Value = 449
$ws.Range("J25").Value = 108
$ws.Range("J27").Value = 131
$ws.Range("J29").Value = 1208
$ws.Range("J31").Value = 196
$ws.Range("J33").Value = 1007
$ws.Range("J36").Value = 296
$ws.Range("J37").Value = 667
$ws.Range("J39").Value = 11
$ws.Range("J42").Value = 908
$ws.Range("J43").Value = 178
$ws.Range("J47").Value = 165
$ws.Range("J48").Value = 258
$ws.Range("J51").Value = 268
$ws.Range("J52").Value = 542
$ws.Range("J54").Value = 423
$ws.Range("J55").Value = 294
$ws.Range("J58").Value = 13
$ws.Range("J60").Value = 129
$ws.Range("J63").Value = 74
$ws.Range("J65").Value = 550
$ws.Range("J66").Value = 66
$ws.Range("J67").Value = 820
$ws.Range("J68").Value = 41
$ws.Range("J71").Value = 74
$ws.Range("J73").Value = 208
$ws.Range("J76").Value = 328
$ws.Range("J77").Value = 164
$ws.Range("J78").Value = 269
$ws.Range("J83").Value = 440
$ws.Range("J84").Value = 184
$ws.Range("J85").Value = 900
$ws.Range("J86").Value = 136
$ws.Range("J88").Value = 231
$ws.Range("J89").Value = 290
$ws.Range("J91").Value = 244
$ws.Range("J94").Value = 221
$ws.Range("J95").Value = 320
$ws.Range("J98").Value = 154
$ws.Range("J101").Value = 21711

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 76
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 204
$ws.Range("J6").Value = 218
$ws.Range("J7").Value = 820

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 59
$ws.Range("J7").Value = 184

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 84
$ws.Range("J6").Value = 204
$ws.Range("J7").Value = 423

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 365
$ws.Range("J3").Value = 423
$ws.Range("J7").Value = 1208

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J2").Value = 43
$ws.Range("J6").Value = 128
$ws.Range("J7").Value = 258

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 187
$ws.Range("J6").Value = 242
$ws.Range("J7").Value = 639

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 183
$ws.Range("J7").Value = 328

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 197
$ws.Range("J3").Value = 184
$ws.Range("J6").Value = 469
$ws.Range("J7").Value = 908

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 87
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 66
$ws.Range("J6").Value = 147
$ws.Range("J7").Value = 294

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 100
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 124
$ws.Range("J7").Value = 449

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 296

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 196
$ws.Range("J3").Value = 192
$ws.Range("J7").Value = 639

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 121
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 70
$ws.Range("J3").Value = 57
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 11

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J3").Value = 69
$ws.Range("J4").Value = 23
$ws.Range("J5").Value = 7
$ws.Range("J7").Value = 340

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J2").Value = 28
$ws.Range("J3").Value = 37
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 72
$ws.Range("J7").Value = 208

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 108
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 81
$ws.Range("J7").Value = 290

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J4").Value = 16
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 73
$ws.Range("J7").Value = 136

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 71
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 46
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J6").Value = 264
$ws.Range("J7").Value = 900

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 63
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J6").Value = 217
$ws.Range("J7").Value = 542

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 13
